$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GPS信息表")

# Create two new formatted rows (9,10) by copying the format of the last existing
# data row (row 8, which already carries the correct border/alignment styles) so
# that no brand-new style entries get appended to styles.xml.
$ws.Range("A8:G8").Copy()
$ws.Range("A9:G10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 3..5: 允许空值 column flips from 否 to 是
$ws.Range("F3").Value = "是"
$ws.Range("F4").Value = "是"
$ws.Range("F5").Value = "是"

# Rows 8..10 become the old 高度/速度/运行状态 rows (shifted down by two, with
# the 允许空值 column also flipped to 是).
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "高度"
$ws.Range("C8").Value = "height"
$ws.Range("D8").Value = "double"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "是"
$ws.Range("G8").Value = ""

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "速度"
$ws.Range("C9").Value = "speed"
$ws.Range("D9").Value = "double"
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "是"
$ws.Range("G9").Value = ""

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "运行状态"
$ws.Range("C10").Value = "status"
$ws.Range("D10").Value = "tinyint"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = "是"
$ws.Range("G10").Value = "0表示运行正常"

# Rows 6..7 become the two newly inserted 百度经度/百度维度 fields.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "百度经度"
$ws.Range("C6").Value = "bmap_longitude"
$ws.Range("D6").Value = "double"
$ws.Range("F6").Value = "是"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "百度维度"
$ws.Range("C7").Value = "bmap_latitude"
$ws.Range("D7").Value = "double"
$ws.Range("F7").Value = "是"

# Sheet view: GPS信息表 becomes the active/selected tab, selection lands on F7.
$ws.Activate() | Out-Null
$ws.Range("F7").Select() | Out-Null
